$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.903.55'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '3.411.15'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.46'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.18'
$ws.Range("E6").Value = '  -1.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  +6.08%  '

$ws.Range("E9").Value = '  +5.58%  '

$ws.Range("E10").Value = '  +2.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.73'
$ws.Range("E11").Value = '  +1.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000217'
$ws.Range("E12").Value = '  +44.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.16'
$ws.Range("E13").Value = '  +8.94%  '

$ws.Range("E14").Value = '  -0.25%  '

$ws.Range("D15").Value = '3.955.83'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.21'
$ws.Range("E16").Value = '  +6.81%  '

$ws.Range("D17").Value = '3.399.94'
$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.50'
$ws.Range("E18").Value = '  +8.00%  '

$ws.Range("E19").Value = '  +6.66%  '

$ws.Range("D20").Value = '61.895.01'
$ws.Range("E20").Value = '  -0.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '452.40'
$ws.Range("E21").Value = '  +44.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.22'
$ws.Range("E22").Value = '  +7.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.20'
$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("E25").Value = '  +3.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.31'
$ws.Range("E26").Value = '  +14.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.00'
$ws.Range("E27").Value = '  +11.08%  '

$ws.Range("E28").Value = '  +0.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.60'
$ws.Range("E29").Value = '  -2.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.70'
$ws.Range("E30").Value = '  -1.35%  '

$ws.Range("E31").Value = '  +5.67%  '

$ws.Range("E32").Value = '  -0.97%  '

$ws.Range("E33").Value = '  -0.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.76'
$ws.Range("E34").Value = '  -4.70%  '

$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0501'
$ws.Range("E36").Value = '  +3.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.83'
$ws.Range("E37").Value = '  +4.02%  '

$ws.Range("E39").Value = '  +2.38%  '

$ws.Range("E40").Value = '  +6.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.95'
$ws.Range("E41").Value = '  -0.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.317'
$ws.Range("E42").Value = '  -1.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.31'
$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("E44").Value = '  +8.16%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.00'
$ws.Range("E45").Value = '  +1.04%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.54'
$ws.Range("E46").Value = '  +14.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.56'
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.26'
$ws.Range("E48").Value = '  +5.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.144'
$ws.Range("E49").Value = '  +20.01%  '

$ws.Range("E50").Value = '  +8.36%  '

$ws.Range("D51").Value = '3.759.49'
$ws.Range("E51").Value = '  -0.44%  '
